$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "add CRC public v2 to list" - append the CRC v2.0-public dataset rows
# (9 datasets) below the existing table, mirroring how the other
# cohort/version blocks are laid out (cohort, version, dataset, synapse_id,
# release_date).
$cohort   = "CRC"
$version  = "v2.0-public"
$datasets = @(
    "cancer_level_dataset_index",
    "cancer_level_dataset_non_index",
    "cancer_panel_test_level_dataset",
    "imaging_level_dataset",
    "med_onc_note_level_dataset",
    "pathology_report_level_dataset",
    "patient_level_dataset",
    "regimen_cancer_level_dataset",
    "tm_level_dataset"
)
$synIds = @(
    "syn39802294",
    "syn39802300",
    "syn39802305",
    "syn39802310",
    "syn39802316",
    "syn39802321",
    "syn39802324",
    "syn39802332",
    "syn39802339"
)
$dates = @(
    "2022-10",
    "2022-11",
    "2022-12",
    "2022-13",
    "2022-14",
    "2022-15",
    "2022-16",
    "2022-17",
    "2022-18"
)

$startRow = 71
$n = $datasets.Count

for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $cohort
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $version
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $datasets[$i]
}

$ws.Cells.Item($startRow, 4).Value = $synIds[0]
$ws.Cells.Item($startRow, 5).Value = $dates[0]

for ($i = 1; $i -lt $n; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $synIds[$i]
}
for ($i = 1; $i -lt $n; $i++) {
    $ws.Cells.Item($startRow + $i, 5).Value = $dates[$i]
}

# Leave the view scrolled down with the cell below the newly typed data
# selected, matching where the user's cursor ended up after data entry.
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C82").Select() | Out-Null
